# Generate Report for Handback
# Updates the localization-status workbook to reflect a handback event:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - "Latest Target File" / "Latest Handback File" columns (E/F) get populated
#    with the same file references as the source / handoff file columns (A/C)
#  - "Latest Handback DateTime" column (G) gets a real timestamp instead of
#    the zero-date placeholder

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: status text is a shared string also used by the
#     language sheets, update it here too so it stays in sync everywhere.
$ws1.Range("B2").Value = $handedBack
$ws1.Range("C2").Value = $handedBack
$ws1.Range("B3").Value = $handedBack
$ws1.Range("C3").Value = $handedBack

# Hyperlink style used throughout this workbook (underline, cornflower blue)
$linkColor = 15570276

# --- zh-cn sheet ---
$ws2.Range("B2").Value = $handedBack
$ws2.Range("B3").Value = $handedBack

$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0913eca420eb03d8f807fbd36c3c13a41745ed77/e2e/569d02c7-ed69-4da3-bcea-4a677ba8dd86.md", "", "", "569d02c7-ed69-4da3-bcea-4a677ba8dd86.md")
$ws2.Range("E2").Font.Underline = $true
$ws2.Range("E2").Font.Color = $linkColor

$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a482334c2ad114a56d55f1f2e87dfafa8732039/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/569d02c7-ed69-4da3-bcea-4a677ba8dd86.42c7ec43a3817ca3d7c3df10c708af6be1f05ace.zh-cn.xlf", "", "", "569d02c7-ed69-4da3-bcea-4a677ba8dd86.42c7ec43a3817ca3d7c3df10c708af6be1f05ace.zh-cn.xlf")
$ws2.Range("F2").Font.Underline = $true
$ws2.Range("F2").Font.Color = $linkColor

$ws2.Range("G2").Value = "2016-03-03 07:40:01"

$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0913eca420eb03d8f807fbd36c3c13a41745ed77/e2e/075dbb6d-8823-4971-95d4-ae348f36376f.md", "", "", "075dbb6d-8823-4971-95d4-ae348f36376f.md")
$ws2.Range("E3").Font.Underline = $true
$ws2.Range("E3").Font.Color = $linkColor

$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a482334c2ad114a56d55f1f2e87dfafa8732039/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/075dbb6d-8823-4971-95d4-ae348f36376f.54c029fd999f89048d96c70d10824dcf335cd76f.zh-cn.xlf", "", "", "075dbb6d-8823-4971-95d4-ae348f36376f.54c029fd999f89048d96c70d10824dcf335cd76f.zh-cn.xlf")
$ws2.Range("F3").Font.Underline = $true
$ws2.Range("F3").Font.Color = $linkColor

$ws2.Range("G3").Value = "2016-03-03 07:40:01"

# --- de-de sheet ---
$ws3.Range("B2").Value = $handedBack
$ws3.Range("B3").Value = $handedBack

$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0913eca420eb03d8f807fbd36c3c13a41745ed77/e2e/569d02c7-ed69-4da3-bcea-4a677ba8dd86.md", "", "", "569d02c7-ed69-4da3-bcea-4a677ba8dd86.md")
$ws3.Range("E2").Font.Underline = $true
$ws3.Range("E2").Font.Color = $linkColor

$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94649cab9cca7f1dea803f5a16e2b0d137102be0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/569d02c7-ed69-4da3-bcea-4a677ba8dd86.42c7ec43a3817ca3d7c3df10c708af6be1f05ace.de-de.xlf", "", "", "569d02c7-ed69-4da3-bcea-4a677ba8dd86.42c7ec43a3817ca3d7c3df10c708af6be1f05ace.de-de.xlf")
$ws3.Range("F2").Font.Underline = $true
$ws3.Range("F2").Font.Color = $linkColor

$ws3.Range("G2").Value = "2016-03-03 07:40:21"

$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0913eca420eb03d8f807fbd36c3c13a41745ed77/e2e/075dbb6d-8823-4971-95d4-ae348f36376f.md", "", "", "075dbb6d-8823-4971-95d4-ae348f36376f.md")
$ws3.Range("E3").Font.Underline = $true
$ws3.Range("E3").Font.Color = $linkColor

$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94649cab9cca7f1dea803f5a16e2b0d137102be0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/075dbb6d-8823-4971-95d4-ae348f36376f.54c029fd999f89048d96c70d10824dcf335cd76f.de-de.xlf", "", "", "075dbb6d-8823-4971-95d4-ae348f36376f.54c029fd999f89048d96c70d10824dcf335cd76f.de-de.xlf")
$ws3.Range("F3").Font.Underline = $true
$ws3.Range("F3").Font.Color = $linkColor

$ws3.Range("G3").Value = "2016-03-03 07:40:21"
